# Non-Onco Study Design functionality
# Merge the "scenario2" (Edit_source_value) rows into the "scenario1" rows as new
# columns, split the "Source_Template" column into "Source_Template_invalid" /
# "Source_Template_valid", and drop the now-empty scenario2 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Pull the "Edit_source_value" data (currently rows 7-10, columns A-F) into
#    memory before we start reshaping the grid.
# ---------------------------------------------------------------------------
$editSourceName   = $ws.Cells.Item(7, 3).Value2     # AutomationTest
$editUpdateName   = $ws.Cells.Item(7, 4).Value2     # Automation_Test_Update
$editAbbrev       = $ws.Cells.Item(8, 3).Value2     # AUT
$editUpdateAbbrev = $ws.Cells.Item(8, 4).Value2     # AUT_UPDT
$editStartOld     = $ws.Cells.Item(9, 3).Value2     # 01/01/2023
$editStartNew     = $ws.Cells.Item(9, 4).Value2     # 01/03/2023
$editEndOld       = $ws.Cells.Item(10, 3).Value2    # 01/05/2023
$editEndNew       = $ws.Cells.Item(10, 4).Value2    # 01/13/2023
$newTemplatePath  = $ws.Cells.Item(7, 6).Value2     # ...AUT_2022_Template - new.xlsx

# ---------------------------------------------------------------------------
# 2) Drop the scenario2 block (rows 7-10) entirely - its data now lives on the
#    scenario1 rows as extra columns.
# ---------------------------------------------------------------------------
$ws.Range("A7:F10").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 3) Make room for the new "Source_Template_valid" column after the existing
#    "Source_Template" column (F).
# ---------------------------------------------------------------------------
$ws.Range("G1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 4) Fill in column D ("Edit_source_value") for the 4 scenario1 rows, copying
#    the quote-prefixed date style from C4/C5 so the two date rows keep the
#    text-as-date formatting instead of becoming real serial dates.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 4).Value = $editUpdateName
$ws.Cells.Item(3, 4).Value = $editUpdateAbbrev

$ws.Cells.Item(4, 3).Copy()
$ws.Cells.Item(4, 4).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(4, 4).Value = "'" + $editStartNew

$ws.Cells.Item(5, 3).Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(5, 4).Value = "'" + $editEndNew
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5) Split "Source_Template" (F) into "Source_Template_invalid" (F, keeps the
#    old path) and "Source_Template_valid" (G, the new path).
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 6).Value = "Source_Template_invalid"
$ws.Cells.Item(1, 7).Value = "Source_Template_valid"
$ws.Cells.Item(2, 7).Value = $newTemplatePath

# Give the new header cell (G1) the same style as the other headers.
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 7).PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 6) Re-create the two placeholder rows (9 and 10) that still carry the
#    (now unused) quote-prefixed date style on C/D but no values.
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 3).Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("D10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 7) Column widths: split the old combined C:D band, give D its own (wider)
#    width, and extend the template-path width onto the new column G.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 21.67
$ws.Columns.Item(7).ColumnWidth = 59.33

# ---------------------------------------------------------------------------
# 8) Selection / view state.
# ---------------------------------------------------------------------------
$ws.Range("A7:F10").Select()
